# "Path to Graduation" planner update:
#  - Fall 2022 block (rows 4-11) gains a 3rd course (E4/F4) and a 2nd
#    course row (row 5 across all three columns) plus a 3rd Spring course
#    (row 6, columns C/D only).
#  - Fall 2023 block (rows 13-20) is reshuffled: a new Fall course is
#    inserted at row 13 (CYBR 2160), the old row-13 Fall course moves to
#    row 14 (CPSC 4115), a 3rd course row (15) is added for Fall/Spring,
#    and a Summer course (CPSC 6985) is added at E13/F13.
#  - The old Fall 2024 course rows (22-23) are emptied out -- that block
#    now has no courses listed, just the header/Total rows.
#  - The Fall 2025 and Fall 2026 blocks (rows 30-47) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fall 2022 / Spring 2022 / Summer 2022 block (rows 4-6) ---
$ws.Range("E4").Value = "CPSC 2108"
$ws.Range("F4").Value = 3

$ws.Range("A5").Value = "CPSC 4111"
$ws.Range("B5").Value = 3
$ws.Range("E5").Value = "CPSC 6180"
$ws.Range("F5").Value = 3

$ws.Range("C6").Value = "CYBR 3115"
$ws.Range("D6").Value = 3

# --- Fall 2023 / Spring 2023 / Summer 2023 block (rows 13-15) ---
$ws.Range("A13").Value = "CYBR 2160"
$ws.Range("E13").Value = "CPSC 6985"
$ws.Range("F13").Value = 4

$ws.Range("A14").Value = "CPSC 4115"

$ws.Range("A15").Value = "CPSC 6185"
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = "CYBR 3119"
$ws.Range("D15").Value = 3

# --- Fall 2024 / Spring 2024 / Summer 2024 block loses its course rows ---
$ws.Range("A22:D23").ClearContents()

# --- Fall 2025 and Fall 2026 blocks (rows 30-47) are removed entirely ---
$ws.Range("A30:F47").ClearContents()
